# Apply the flowering_success_2021.docx edit:
#  1. Shrink the header-row height of Table 3 and Table 5 from 637 twips
#     (31.85 pt) to 571 twips (28.55 pt) — hRule stays "auto".
#  2. Fix the mojibake "Ï‡" (UTF-8 chi re-decoded as Latin-1) back to the
#     real Greek chi character "χ" in both of those same tables'
#     chi-squared (χ²) column headers.

$d = $word.ActiveDocument

# --- 1. Header row heights -------------------------------------------------
# 637 twips == 31.85 pt ; 571 twips == 28.55 pt (Row.Height is in points).
$oldHeightPt = 637 / 20
$newHeightPt = 571 / 20

for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    $headerRow = $tbl.Rows.Item(1)
    if ([Math]::Abs($headerRow.Height - $oldHeightPt) -lt 0.01) {
        $headerRow.Height = $newHeightPt
    }
}

# --- 2. Mojibake chi character fix -----------------------------------------
$mojibake = "Ï‡"
$fixed = "χ"

$d.Content.Find.Execute(
    $mojibake, $true, $false, $false, $false, $false,
    $true, 1, $false, $fixed, 2
) | Out-Null
